$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift years/values up by one row (drop old 2020 row, add 2025 at the end)
$ws.Range("A2").Value = 2021
$ws.Range("B2").Value = 950468.77

$ws.Range("A3").Value = 2022
$ws.Range("B3").Value = 1766113.68
$ws.Range("C3").Value = 85.81501420609537

$ws.Range("A4").Value = 2023
$ws.Range("B4").Value = 2842827.64
$ws.Range("C4").Value = 60.96515599154411

$ws.Range("A5").Value = 2024
$ws.Range("B5").Value = 4442894.22
$ws.Range("C5").Value = 56.28433315781323

$ws.Range("A6").Value = 2025
$ws.Range("B6").Value = 1014612.85
$ws.Range("C6").Value = -77.1632454035784

# Remove the now-obsolete row 7 entirely
$ws.Rows("7:7").Delete()
